$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting all existing price-history rows
# (previously rows 2..74) down by one (to 3..75).
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest date and the same
# price figures that the rest of the series carries. The date column is
# stored as plain text throughout the sheet ("yyyy-mm-dd"), so force a
# Text number format before writing the value to stop Excel's COM layer
# from auto-converting the literal into a date serial, then restore the
# default ("Normal") style so the cell matches its siblings exactly.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-02"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
